$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("월요일")

# Clear A8 entirely (value and remains empty)
$ws.Range("A8").ClearContents()

# Clear B8 and C8 contents but keep their number-format style
$ws.Range("B8:C8").ClearContents()

# Update the active selection to A8
$ws.Range("A8").Select()
